# Insert a new row before row 883 ("Curso Preparatório de Física" / "António Casaca")
# for the new "Academia MikroTik" link, shifting every row from 883 onward down by one
# (old A1:E1141 -> new A1:E1142).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(883).EntireRow.Insert()

$ws.Cells.Item(883, 1).Value = "Cursos"
$ws.Cells.Item(883, 2).Value = "Outros Cursos"
$ws.Cells.Item(883, 3).Value = "Academia MikroTik"
$ws.Cells.Item(883, 4).Value = "aqui"
$ws.Cells.Item(883, 5).Value = "https://www.isel.pt/sites/default/files/002_pdf/50_MikroTikAcademy_EN.pdf"
